$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.385.64'
$ws.Range("E2").Value = '  +0.85%  '
$ws.Range("D3").Value = '1.669.84'
$ws.Range("E3").Value = '  +1.03%  '
$ws.Range("E4").Value = '  +0.56%  '
$ws.Range("D5").Value = '''221.11'
$ws.Range("E5").Value = '  +1.63%  '
$ws.Range("D6").Value = '''0.5326'
$ws.Range("E6").Value = '  +0.53%  '
$ws.Range("E7").Value = '  +0.50%  '
$ws.Range("D8").Value = '''0.2659'
$ws.Range("E8").Value = '  +1.46%  '
$ws.Range("D9").Value = '''0.06369'
$ws.Range("E9").Value = '  +0.75%  '
$ws.Range("D10").Value = '''20.83'
$ws.Range("E10").Value = '  +2.22%  '
$ws.Range("D11").Value = '''0.07855'
$ws.Range("E11").Value = '  +0.56%  '
$ws.Range("D12").Value = '''4.524'
$ws.Range("E12").Value = '  +0.19%  '
$ws.Range("D13").Value = '1.673.07'
$ws.Range("E13").Value = '  +0.76%  '
$ws.Range("D14").Value = '1.900.15'
$ws.Range("E14").Value = '  +1.04%  '
$ws.Range("D15").Value = '''0.5598'
$ws.Range("E15").Value = '  +2.07%  '
$ws.Range("D16").Value = '0.0₅8183'
$ws.Range("E16").Value = '  +0.08%  '
$ws.Range("D17").Value = '''66.12'
$ws.Range("E17").Value = '  +1.20%  '
$ws.Range("D18").Value = '26.410.46'
$ws.Range("E18").Value = '  +0.98%  '
$ws.Range("E19").Value = '  +0.57%  '
$ws.Range("D20").Value = '''4.713'
$ws.Range("E20").Value = '  +2.59%  '
$ws.Range("D21").Value = '''197.77'
$ws.Range("E21").Value = '  +3.54%  '
$ws.Range("D22").Value = '''10.29'
$ws.Range("E22").Value = '  +2.11%  '
$ws.Range("D23").Value = '''6.064'
$ws.Range("E23").Value = '  +1.03%  '
$ws.Range("D24").Value = '''1.012'
$ws.Range("E24").Value = '  +0.47%  '
$ws.Range("D25").Value = '''145.63'
$ws.Range("E25").Value = '  +0.26%  '
$ws.Range("D26").Value = '''0.1224'
$ws.Range("E26").Value = '  +0.02%  '
$ws.Range("D27").Value = '''7.246'
$ws.Range("E27").Value = '  +0.72%  '
$ws.Range("D28").Value = '''16.18'
$ws.Range("E28").Value = '  +1.75%  '
$ws.Range("E29").Value = '  +2.66%  '
$ws.Range("D30").Value = '''0.05914'
$ws.Range("E30").Value = '  +3.43%  '
$ws.Range("E31").Value = '  +1.20%  '
$ws.Range("D32").Value = '''3.561'
$ws.Range("E32").Value = '  +0.47%  '
$ws.Range("D33").Value = '''3.327'
$ws.Range("E33").Value = '  +2.02%  '
$ws.Range("D34").Value = '''1.607'
$ws.Range("E34").Value = '  +1.20%  '
$ws.Range("D35").Value = '''0.9676'
$ws.Range("E35").Value = '  +2.10%  '
$ws.Range("D36").Value = '''2.839'
$ws.Range("E36").Value = '  +1.24%  '
$ws.Range("D37").Value = '''2.439'
$ws.Range("E37").Value = '  +0.68%  '
$ws.Range("D38").Value = '''0.5832'
$ws.Range("E38").Value = '  +2.08%  '
$ws.Range("D39").Value = '''0.01616'
$ws.Range("E39").Value = '  +0.60%  '
$ws.Range("D40").Value = '1.078.13'
$ws.Range("E40").Value = '  +3.83%  '
$ws.Range("D41").Value = '''5.933'
$ws.Range("E41").Value = '  +2.41%  '
$ws.Range("D42").Value = '''0.8639'
$ws.Range("E42").Value = '  +1.57%  '
$ws.Range("D44").Value = '''103.03'
$ws.Range("E44").Value = '  -0.74%  '
$ws.Range("D45").Value = '1.811.00'
$ws.Range("E45").Value = '  +0.98%  '
$ws.Range("D46").Value = '''58.45'
$ws.Range("E46").Value = '  +3.10%  '
$ws.Range("E47").Value = '  +5.84%  '
$ws.Range("E48").Value = '  +0.76%  '
$ws.Range("E49").Value = '  +1.45%  '
$ws.Range("D50").Value = '''7.990'
$ws.Range("E50").Value = '  +1.93%  '
$ws.Range("D51").Value = '''0.05158'
$ws.Range("E51").Value = '  +0.06%  '
